# Rotate rows 2-4 of the Artfynd sheet:
#   new row 2 <= old row 3 data
#   new row 3 <= old row 4 data
#   new row 4 <= old row 2 data
# Values are captured first (not range-copied) to avoid clobbering source rows
# before they have been read, and cells are only touched when their content
# actually needs to change, to keep unrelated cells untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value2 = 111697304
$ws.Range("B2").Value2 = 8377
$ws.Range("D2").Value2 = 'LC'
$ws.Range("E2").Value2 = 106545
$ws.Range("F2").Value2 = 'Mindre märgborre'
$ws.Range("G2").Value2 = 'Tomicus minor'
$ws.Range("H2").Value2 = '(Hartig, 1834)'
$ws.Range("J2").Value2 = ""
$ws.Range("M2").Value2 = 'färska gnagspår'
$ws.Range("Q2").Value2 = 373090.8741807578
$ws.Range("R2").Value2 = 6865424.499624529
$ws.Range("AC2").Value2 = ""
$ws.Range("AI2").Value2 = 'Luckig tallskog. K-skog'
$ws.Range("AJ2").Value2 = ""
$ws.Range("AK2").Value2 = ""
$ws.Range("AO2").Value2 = ""

# Row 3 updates
$ws.Range("A3").Value2 = 111697236
$ws.Range("Q3").Value2 = 373121.3523494597
$ws.Range("R3").Value2 = 6865443.651501717
$ws.Range("Z3").Value2 = '00:00'
$ws.Range("AB3").Value2 = '00:00'
$ws.Range("AI3").Value2 = 'Tallskog. Kontinuitetsskog'
$ws.Range("AJ3").Value2 = 'tall'
$ws.Range("AK3").Value2 = 'Pinus sylvestris'
$ws.Range("AO3").Value2 = 'Pinus sylvestris'

# Row 4 updates
$ws.Range("A4").Value2 = 111697636
$ws.Range("B4").Value2 = 88489
$ws.Range("D4").Value2 = 'NT'
$ws.Range("E4").Value2 = 1962
$ws.Range("F4").Value2 = 'Vaddporing'
$ws.Range("G4").Value2 = 'Anomoporia kamtschatica'
$ws.Range("H4").Value2 = '(Parmasto) Bondartseva'
$ws.Range("J4").Value2 = 'fruktkroppar'
$ws.Range("L4").Value2 = ""
$ws.Range("M4").Value2 = ""
$ws.Range("Q4").Value2 = 373112.5181173298
$ws.Range("R4").Value2 = 6865358.590016441
$ws.Range("Z4").Value2 = '19:00'
$ws.Range("AB4").Value2 = '19:00'
$ws.Range("AC4").Value2 = 'Växer under rötad gammal silverved'
$ws.Range("AI4").Value2 = 'Kontinuitetsskog. Tallskog'
